$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Datos actualizados" timestamp string in A1
$ws.Range("A1").Value2 = "Datos actualizados a 17 de Junio de 2020 a las 14:58"

# Row 13
$ws.Range("B13").Value2 = 188523
$ws.Range("C13").Value2 = 141
$ws.Range("E13").Value2 = 6013

# Row 19
$ws.Range("B19").Value2 = 141234
$ws.Range("C19").Value2 = 4919
$ws.Range("D19").Value2 = 91662
$ws.Range("E19").Value2 = 48481
$ws.Range("G19").Value2 = 39
$ws.Range("H19").Value2 = 1091

# Row 23
$ws.Range("B23").Value2 = 83174
$ws.Range("C23").Value2 = 1097
$ws.Range("D23").Value2 = 62172
$ws.Range("E23").Value2 = 20920
$ws.Range("G23").Value2 = 2
$ws.Range("H23").Value2 = 82

# Row 29
$ws.Range("B29").Value2 = 49204
$ws.Range("C29").Value2 = 117
$ws.Range("G29").Value2 = 4
$ws.Range("H29").Value2 = 6074

# Row 35
$ws.Range("A35").Value2 = "Portugal"
$ws.Range("B35").Value2 = 37672
$ws.Range("C35").Value2 = 336
$ws.Range("D35").Value2 = 23580
$ws.Range("E35").Value2 = 12569
$ws.Range("G35").Value2 = 1
$ws.Range("H35").Value2 = 1523

# Row 36
$ws.Range("A36").Value2 = "Kuwait"
$ws.Range("B36").Value2 = 37533
$ws.Range("C36").Value2 = 575
$ws.Range("D36").Value2 = 28896
$ws.Range("E36").Value2 = 8331
$ws.Range("G36").Value2 = 3
$ws.Range("H36").Value2 = 306

# Row 37
$ws.Range("D37").Value2 = 10512
$ws.Range("E37").Value2 = 22761
$ws.Range("G37").Value2 = 8
$ws.Range("H37").Value2 = 886

# Row 51
$ws.Range("E51").Value2 = 5639
$ws.Range("G51").Value2 = 1
$ws.Range("H51").Value2 = 48

# Row 58
$ws.Range("A58").Value2 = "Dinamarca"
$ws.Range("B58").Value2 = 12294
$ws.Range("C58").Value2 = 44
$ws.Range("D58").Value2 = 11185
$ws.Range("E58").Value2 = 511
$ws.Range("G58").Value2 = 0
$ws.Range("H58").Value2 = 598

# Row 59
$ws.Range("A59").Value2 = "Moldavia"
$ws.Range("B59").Value2 = 12254
$ws.Range("D59").Value2 = 7077
$ws.Range("E59").Value2 = 4750
$ws.Range("G59").Value2 = 4
$ws.Range("H59").Value2 = 427

# Row 69
$ws.Range("E69").Value2 = 279
$ws.Range("G69").Value2 = 1
$ws.Range("H69").Value2 = 243

# Row 76
$ws.Range("B76").Value2 = 5638
$ws.Range("C76").Value2 = 145
$ws.Range("E76").Value2 = 1523

# Row 83
$ws.Range("B83").Value2 = 4482
$ws.Range("C83").Value2 = 183
$ws.Range("D83").Value2 = 1803
$ws.Range("E83").Value2 = 2469
$ws.Range("G83").Value2 = 9
$ws.Range("H83").Value2 = 210

# Row 89
$ws.Range("B89").Value2 = 3759
$ws.Range("C89").Value2 = 129
$ws.Range("D89").Value2 = 849
$ws.Range("E89").Value2 = 2847
$ws.Range("G89").Value2 = 2
$ws.Range("H89").Value2 = 63

# Row 93
$ws.Range("A93").Value2 = "Bosnia y Herzegovina"
$ws.Range("B93").Value2 = 3141
$ws.Range("C93").Value2 = 56
$ws.Range("D93").Value2 = 2197
$ws.Range("E93").Value2 = 776
$ws.Range("H93").Value2 = 168

# Row 94
$ws.Range("A94").Value2 = "Tailandia"
$ws.Range("B94").Value2 = 3135
$ws.Range("D94").Value2 = 2996
$ws.Range("E94").Value2 = 81
$ws.Range("H94").Value2 = 58

# Row 100
$ws.Range("B100").Value2 = 2258
$ws.Range("C100").Value2 = 3
$ws.Range("D100").Value2 = 2141
$ws.Range("E100").Value2 = 10

# Row 104
$ws.Range("B104").Value2 = 1921
$ws.Range("C104").Value2 = 6
$ws.Range("E104").Value2 = 513

# Row 168
$ws.Range("D168").Value2 = 176
$ws.Range("E168").Value2 = 0

# Row 206
$ws.Range("A206").Value2 = "Groenlandia"

# Row 207
$ws.Range("A207").Value2 = "Islas Malvinas"

# Row 210
$ws.Range("A210").Value2 = "Seychelles"
$ws.Range("D210").Value2 = 11
$ws.Range("H210").Value2 = 0

# Row 211
$ws.Range("A211").Value2 = "Montserrat"
$ws.Range("D211").Value2 = 10
$ws.Range("H211").Value2 = 1

# Row 213
$ws.Range("A213").Value2 = "Papua Nueva Guinea"
$ws.Range("D213").Value2 = 8
$ws.Range("H213").Value2 = 0

# Row 214
$ws.Range("A214").Value2 = "Islas Virgenes Britanicas"
$ws.Range("D214").Value2 = 7
$ws.Range("H214").Value2 = 1
